$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Re-center the previous "record" block (rows 19-24) --------------------
# These cells currently carry style slot 2 (center/center). Nudging the
# alignment (no-op value-wise) lets the engine settle them onto the
# equivalent pre-existing style slot 1, matching the target file.
$ws.Range("E19:E24").HorizontalAlignment = -4108
$ws.Range("E19:E24").VerticalAlignment = -4108
$ws.Range("A20:C24").HorizontalAlignment = -4108
$ws.Range("A20:C24").VerticalAlignment = -4108

# --- Append the new diary entry (rows 25-27) --------------------------------
$ws.Range("A25").Value = "今天的话，学校要考英语，给做了一下，然后乡政府领钱那个要我准备一份稿子，就去网上找了然后改了一段时间"
$ws.Range("E25").Value = "今天把第三章给过了，循环啊，条件语句之类的，那上面的例子，之前全写过，就偷了一波懒"
$ws.Range("E27").Value = "然后今天的那个代码，就是我把之前那个运营不了的代码给又抄了一遍"

# Merge the new block the same way the previous ones were merged.
$ws.Range("A25:C27").Merge()
$ws.Range("E25:E26").Merge()

# Center the new block (A25:C27 holds the date/status text, E25:E26 holds
# the first journal paragraph; E27 -- a continuation note -- stays default).
$ws.Range("A25:C27").HorizontalAlignment = -4108
$ws.Range("A25:C27").VerticalAlignment = -4108
$ws.Range("E25:E26").HorizontalAlignment = -4108
$ws.Range("E25:E26").VerticalAlignment = -4108

# Match the author's final selection/cursor position.
$ws.Range("A25:C27").Select()
